$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are written as text, matching the source data
# (e.g. "25.750.22", "1.003") instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.750.22"
$ws.Range("E2").Value = "  +5.28%  "

$ws.Range("D3").Value = "1.702.96"
$ws.Range("E3").Value = "  +3.22%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "330.43"
$ws.Range("E5").Value = "  +5.95%  "

$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "0.3681"
$ws.Range("E7").Value = "  +0.85%  "

$ws.Range("D8").Value = "48.29"
$ws.Range("E8").Value = "  +3.66%  "

$ws.Range("D9").Value = "0.3304"
$ws.Range("E9").Value = "  +1.75%  "

$ws.Range("D10").Value = "1.166"
$ws.Range("E10").Value = "  +3.62%  "

$ws.Range("D11").Value = "0.07333"
$ws.Range("E11").Value = "  +4.25%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "6.183"
$ws.Range("E13").Value = "  +3.64%  "

$ws.Range("D14").Value = "19.99"
$ws.Range("E14").Value = "  +3.02%  "

$ws.Range("D15").Value = "6.851"
$ws.Range("E15").Value = "  +3.65%  "

$ws.Range("D16").Value = "1.701.56"
$ws.Range("E16").Value = "  +3.14%  "

$ws.Range("D17").Value = "0.00001064"
$ws.Range("E17").Value = "  +2.25%  "

$ws.Range("D18").Value = "0.06620"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("D19").Value = "80.98"
$ws.Range("E19").Value = "  +2.94%  "

$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("D21").Value = "16.13"
$ws.Range("E21").Value = "  +3.01%  "

$ws.Range("D22").Value = "6.033"
$ws.Range("E22").Value = "  +1.70%  "

$ws.Range("D23").Value = "12.93"
$ws.Range("E23").Value = "  +3.16%  "

$ws.Range("D24").Value = "25.737.43"
$ws.Range("E24").Value = "  +5.36%  "

$ws.Range("D25").Value = "2.457"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").Value = "2.475"
$ws.Range("E26").Value = "  +6.26%  "

$ws.Range("E27").Value = "  +1.92%  "

$ws.Range("D28").Value = "19.13"
$ws.Range("E28").Value = "  +2.90%  "

$ws.Range("D29").Value = "1.293"
$ws.Range("E29").Value = "  +8.73%  "

$ws.Range("D30").Value = "1.892.25"
$ws.Range("E30").Value = "  +3.36%  "

$ws.Range("D31").Value = "127.80"
$ws.Range("E31").Value = "  +3.29%  "

$ws.Range("D32").Value = "4.116"
$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("D33").Value = "5.914"
$ws.Range("E33").Value = "  +3.77%  "

$ws.Range("D34").Value = "0.08484"
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("D35").Value = "1.704"
$ws.Range("E35").Value = "  +3.27%  "

$ws.Range("D36").Value = "12.80"
$ws.Range("E36").Value = "  +5.59%  "

$ws.Range("D37").Value = "5.309"
$ws.Range("E37").Value = "  +1.74%  "

$ws.Range("D38").Value = "1.272"
$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("D39").Value = "0.06193"
$ws.Range("E39").Value = "  +2.81%  "

$ws.Range("D40").Value = "8.494"
$ws.Range("E40").Value = "  +4.18%  "

$ws.Range("D41").Value = "0.2113"
$ws.Range("E41").Value = "  +2.55%  "

$ws.Range("D42").Value = "0.02242"
$ws.Range("E42").Value = "  +0.56%  "

$ws.Range("D43").Value = "14.75"
$ws.Range("E43").Value = "  +16.90%  "

$ws.Range("D44").Value = "0.6095"
$ws.Range("E44").Value = "  +3.26%  "

$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.20%  "

$ws.Range("D46").Value = "3.836"
$ws.Range("E46").Value = "  +1.96%  "

$ws.Range("D47").Value = "0.5814"
$ws.Range("E47").Value = "  +3.48%  "

$ws.Range("D48").Value = "126.27"
$ws.Range("E48").Value = "  +2.66%  "

$ws.Range("D49").Value = "1.997"
$ws.Range("E49").Value = "  +2.53%  "

$ws.Range("D50").Value = "0.07215"
$ws.Range("E50").Value = "  +4.51%  "

$ws.Range("D51").Value = "76.33"
$ws.Range("E51").Value = "  +2.68%  "
